$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.792.91"
$ws.Range("E2").Value = "  +7.25%  "

$ws.Range("D3").Value = "1.950.97"
$ws.Range("E3").Value = "  +5.53%  "

$ws.Range("E4").Value = "  -0.51%  "

$ws.Range("D5").Value = "'341.79"
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").Value = "'0.4787"
$ws.Range("E7").Value = "  +3.04%  "

$ws.Range("D8").Value = "'0.4137"
$ws.Range("E8").Value = "  +7.04%  "

$ws.Range("D9").Value = "'47.85"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").Value = "'0.08244"
$ws.Range("E10").Value = "  +4.35%  "

$ws.Range("D11").Value = "'1.039"
$ws.Range("E11").Value = "  +7.31%  "

$ws.Range("D12").Value = "'22.81"
$ws.Range("E12").Value = "  +7.09%  "

$ws.Range("D13").Value = "1.955.43"
$ws.Range("E13").Value = "  +7.01%  "

$ws.Range("E14").Value = "  +4.62%  "

$ws.Range("D15").Value = "'7.392"
$ws.Range("E15").Value = "  +3.26%  "

$ws.Range("D16").Value = "'92.01"
$ws.Range("E16").Value = "  +2.16%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "'0.00001060"
$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("D19").Value = "'0.06694"
$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("D20").Value = "'18.05"
$ws.Range("E20").Value = "  +4.03%  "

$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "29.746.32"
$ws.Range("E22").Value = "  +7.14%  "

$ws.Range("D23").Value = "'5.602"
$ws.Range("E23").Value = "  +4.84%  "

$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").Value = "2.182.96"

$ws.Range("D27").Value = "'160.85"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("D28").Value = "'20.25"
$ws.Range("E28").Value = "  +3.89%  "

$ws.Range("D29").Value = "'2.174"
$ws.Range("E29").Value = "  +5.15%  "

$ws.Range("D30").Value = "'5.669"
$ws.Range("E30").Value = "  +5.78%  "

$ws.Range("D31").Value = "'122.92"
$ws.Range("E31").Value = "  +3.52%  "

$ws.Range("D32").Value = "'1.007"
$ws.Range("E32").Value = "  +6.53%  "

$ws.Range("D33").Value = "'0.09647"
$ws.Range("E33").Value = "  +2.52%  "

$ws.Range("E34").Value = "  +11.15%  "

$ws.Range("D35").Value = "'3.679"
$ws.Range("E35").Value = "  +2.61%  "

$ws.Range("D36").Value = "'5.509"
$ws.Range("E36").Value = "  +4.66%  "

$ws.Range("D37").Value = "'0.06264"
$ws.Range("E37").Value = "  +4.32%  "

$ws.Range("D38").Value = "'0.02317"
$ws.Range("E38").Value = "  +4.91%  "

$ws.Range("D39").Value = "'8.493"
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("D40").Value = "'1.187"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").Value = "'0.6090"
$ws.Range("E41").Value = "  +4.85%  "

$ws.Range("E42").Value = "  +6.23%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1897"
$ws.Range("E43").Value = "  +3.02%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.401"
$ws.Range("E45").Value = "  +32.07%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.275"
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("D47").Value = "'0.5722"
$ws.Range("E47").Value = "  +5.05%  "

$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("E49").Value = "  +8.43%  "

$ws.Range("D50").Value = "'1.992"
$ws.Range("E50").Value = "  +3.34%  "

$ws.Range("D51").Value = "'113.16"
$ws.Range("E51").Value = "  +2.04%  "
